# Update the two-digit multiplication answer cells in the table.
# Replacements are applied in document order so that the two cells
# that originally shared the same text ("67x59=3953") each receive
# their own distinct replacement.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "11×50=550";   New = "44×34=1496" },
    @{ Old = "18×60=1080";  New = "62×25=1550" },
    @{ Old = "49×87=4263";  New = "92×45=4140" },
    @{ Old = "98×51=4998";  New = "54×36=1944" },
    @{ Old = "38×19=722";   New = "79×55=4345" },
    @{ Old = "57×27=1539";  New = "92×72=6624" },
    @{ Old = "20×44=880";   New = "14×84=1176" },
    @{ Old = "65×21=1365";  New = "34×72=2448" },
    @{ Old = "43×27=1161";  New = "49×14=686" },
    @{ Old = "76×66=5016";  New = "42×85=3570" },
    @{ Old = "67×59=3953";  New = "95×79=7505" },
    @{ Old = "61×24=1464";  New = "73×33=2409" },
    @{ Old = "70×99=6930";  New = "80×93=7440" },
    @{ Old = "55×71=3905";  New = "99×27=2673" },
    @{ Old = "38×37=1406";  New = "23×22=506" },
    @{ Old = "73×86=6278";  New = "59×21=1239" },
    @{ Old = "21×15=315";   New = "59×22=1298" },
    @{ Old = "81×60=4860";  New = "48×83=3984" },
    @{ Old = "41×73=2993";  New = "67×50=3350" },
    @{ Old = "79×53=4187";  New = "33×65=2145" },
    @{ Old = "79×21=1659";  New = "15×75=1125" },
    @{ Old = "37×46=1702";  New = "68×67=4556" },
    @{ Old = "24×20=480";   New = "23×32=736" },
    @{ Old = "67×59=3953";  New = "42×16=672" },
    @{ Old = "43×68=2924";  New = "96×17=1632" }
)

$searchStart = 0

foreach ($pair in $replacements) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute($pair.Old, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $pair.New
        $searchStart = $rng.End
    }
}
